# Reverse the order of the comma-separated "Recorded By" names in column G
# for every data row on the active sheet (row 1 is the header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ",\s*"
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
